$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second applicant's row (was row 6) and the trailing blank row (was row 7)
$ws.Rows("6:7").Delete() | Out-Null

# Update the remaining applicant's card number
$ws.Range("B5").Value = 157

# "Trình độ" value DL3 -> DL1 for the remaining record
$ws.Range("E5").Value = "DL1"

# Fill in the actual Decision number / day / month / year values
$ws.Range("H2").Value = 1233
$ws.Range("F3").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = 2019

# Move the active selection to B5
$ws.Range("B5").Select() | Out-Null
